# Auto commit at 2025-10-15  8:14:01.73
# Append two new daily rows (2025-10-14) for the two charging stations to
# Sheet1, right after the existing last row (87).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 88 - 四方坪站充电量(kw), date serial 45944 (2025-10-14)
$ws.Range("A88").Value = 45944
$ws.Range("B88").Value = "四方坪站充电量(kw)"
$ws.Range("C88").Value = 722.47399999999982
$ws.Range("D88").Value = 1269.9000000000001
$ws.Range("E88").Value = 337.05599999999993
$ws.Range("F88").Value = 403.80700000000002
$ws.Range("G88").Value = 446.70000000000005
$ws.Range("H88").Value = 556.23
$ws.Range("I88").Value = 602.20400000000006
$ws.Range("J88").Value = 205.35599999999999
$ws.Range("K88").Value = 130.86099999999999
$ws.Range("L88").Value = 181.23900000000003
$ws.Range("M88").Value = 158.60200000000003
$ws.Range("N88").Value = 225.41
$ws.Range("O88").Value = 1061.979
$ws.Range("P88").Value = 893.90100000000018
$ws.Range("Q88").Value = 419.84399999999999
$ws.Range("R88").Value = 257.52000000000004
$ws.Range("S88").Value = 306.30499999999995
$ws.Range("T88").Value = 265.41699999999997
$ws.Range("U88").Value = 86.798000000000002
$ws.Range("V88").Value = 143.863
$ws.Range("W88").Value = 24.54
$ws.Range("X88").Value = 21.1
$ws.Range("Y88").Value = 11.602
$ws.Range("Z88").Value = 42.647999999999996

# Row 89 - 高岭站充电量(kw), date serial 45944 (2025-10-14)
$ws.Range("A89").Value = 45944
$ws.Range("B89").Value = "高岭站充电量(kw)"
$ws.Range("C89").Value = 434.911
$ws.Range("D89").Value = 287.24900000000002
$ws.Range("E89").Value = 41.295999999999999
$ws.Range("F89").Value = 111.29299999999998
$ws.Range("G89").Value = 57.765000000000001
$ws.Range("H89").Value = 237.15
$ws.Range("I89").Value = 84.281000000000006
$ws.Range("J89").Value = 223.36600000000001
$ws.Range("K89").Value = 328.72599999999994
$ws.Range("L89").Value = 120.17099999999999
$ws.Range("M89").Value = 82.073999999999998
$ws.Range("N89").Value = 78.890000000000015
$ws.Range("O89").Value = 470.36900000000003
$ws.Range("P89").Value = 517.41399999999987
$ws.Range("Q89").Value = 328.63200000000001
$ws.Range("R89").Value = 156.256
$ws.Range("S89").Value = 132.08100000000002
$ws.Range("T89").Value = 121.26499999999999
$ws.Range("U89").Value = 96.436999999999998
$ws.Range("V89").Value = 70.354000000000013
$ws.Range("W89").Value = 18.972000000000001
$ws.Range("X89").Value = 8.9529999999999994
$ws.Range("Y89").Value = 22.055
$ws.Range("Z89").Value = 26.774000000000001

# Match the author's final on-screen selection after the edit.
[void]$ws.Range("F97").Select()
